# "switched left & right again"
# The Plane column (C) values for mic groups B1-B17 and D1-D17 were swapped:
#  - B1..B17 (rows 2-18)  : "left"  -> "right"
#  - D1..D17 (rows 36-52) : "right" -> "left"
# Also restore the scrolled viewport (topLeftCell = A39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 18; $r++) {
    $ws.Range("C$r").Value = "right"
}

for ($r = 36; $r -le 52; $r++) {
    $ws.Range("C$r").Value = "left"
}

# Scroll the sheet view so row 39 is at the top-left of the visible area.
$ws.Application.ActiveWindow.ScrollRow = 39
